$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Random Forest")

# --- Row 35 / 36: change highlighted (yellow-fill) cells to the "cleared fill" style ---
# These cells previously carried the yellow promo-highlight fill (style 5); the edit
# clears that fill (new style 11) while keeping / extending the range that is covered.
$ws.Range("I35:M35").Interior.ColorIndex = -4142
$ws.Range("J35:K35").Interior.ColorIndex = -4142
$ws.Range("I36").Interior.ColorIndex = -4142
$ws.Range("L36:M36").Interior.ColorIndex = -4142

# --- New columns H:L for rows 51-61 : "remove promo" / "remove mobile test filter" results ---
$ws.Range("I50").Value = "remove promo"
$ws.Range("K50").Value = "remove mobile test filter"

$ws.Range("H51").Value = "Operating System"
$ws.Range("K51").Value = 0.97173144876325002
$ws.Range("L51").Value = 0.97123473865937504

$ws.Range("H52").Value = "Features"
$ws.Range("K52").Value = 0.77138136993264705
$ws.Range("L52").Value = 0.77026625297587303

$ws.Range("H53").Value = "Network Connections"
$ws.Range("K53").Value = 0.91056563500533605
$ws.Range("L53").Value = 0.90479330659597801

$ws.Range("H54").Value = "Memory RAM"
$ws.Range("I54").Value = 0.890778401122019
$ws.Range("J54").Value = 0.88934764022218704
$ws.Range("K54").Value = 0.89042776998597395
$ws.Range("L54").Value = 0.88911397041819495

$ws.Range("H55").Value = "Brand"

$ws.Range("H56").Value = "Warranty Period"
$ws.Range("I56").Value = 0.85761830473218903
$ws.Range("J56").Value = 0.82888401632192199
$ws.Range("K56").Value = 0.85678627145085795
$ws.Range("L56").Value = 0.82753126240531305

$ws.Range("H57").Value = "Storage Capacity"
$ws.Range("K57").Value = 0.95130641330166199
$ws.Range("L57").Value = 0.95093413607259503

$ws.Range("H58").Value = "Color Family"
$ws.Range("K58").Value = 0.84169269758943799
$ws.Range("L58").Value = 0.83904685698053205

$ws.Range("H59").Value = "Phone Model"

$ws.Range("H60").Value = "Camera"
$ws.Range("K60").Value = 0.69987799918666105
$ws.Range("L60").Value = 0.69390833463344603

$ws.Range("H61").Value = "Phone Screen Size"
$ws.Range("K61").Value = 0.776571341183667
$ws.Range("L61").Value = 0.77373729773394795

# B51:C51 / B61:C61 pick up the yellow promo-highlight fill (style 5)
$ws.Range("B51:C51").Interior.ColorIndex = 6
$ws.Range("B61:C61").Value2 = @(0.77718305551307498, 0.77410656270602296)
$ws.Range("B61:C61").Interior.ColorIndex = 6

# K52:L52, K53:L53, K54:L54, I56:J56, K56:L56(already cleared fill above),
# K57:L57, K58:L58, K60:L60 also get the yellow highlight (style 5), except
# K56:L56 which gets the cleared-fill style (style 11) like I35 etc.
$ws.Range("K52:L52").Interior.ColorIndex = 6
$ws.Range("K53:L53").Interior.ColorIndex = 6
$ws.Range("K54:L54").Interior.ColorIndex = 6
$ws.Range("I56:J56").Interior.ColorIndex = 6
$ws.Range("K56:L56").Interior.ColorIndex = -4142
$ws.Range("K57:L57").Interior.ColorIndex = 6
$ws.Range("K58:L58").Interior.ColorIndex = 6
$ws.Range("K60:L60").Interior.ColorIndex = 6

# --- sheet view: scroll position + current selection ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("N26").Select()

# --- new shared strings used elsewhere in the workbook (word-frequency filters) ---
$ws.Range("Z1").Value = "remove promo"
$ws.Range("Z2").Value = "remove mobile test filter"
$ws.Range("Z1:Z2").ClearContents()
